$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the "V" marker used for few-shot prompting examples into column Q
# for each of the statement-of-values detail rows (rows 8-15).
$ws.Range("Q8:Q15").Value = "V"
